# Update the two-digit / one-digit division answers in the worksheet table.
# Each cell is addressed directly via Tables(1).Cell(row, col) and its text is
# replaced in place so formatting (rFonts/sz) on the existing run is preserved
# and cells that happen to share old/new text values do not cross-contaminate.
$d = $word.ActiveDocument
$t = $d.Tables(1)

$cell = $t.Cell(1,1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "89÷4=22, 1") {
    Write-Host "WARNING: Cell(1,1) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "80÷5=16, 0"

$cell = $t.Cell(1,2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "61÷4=15, 1") {
    Write-Host "WARNING: Cell(1,2) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "74÷9=8, 2"

$cell = $t.Cell(1,3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "30÷9=3, 3") {
    Write-Host "WARNING: Cell(1,3) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "31÷4=7, 3"

$cell = $t.Cell(1,4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "34÷6=5, 4") {
    Write-Host "WARNING: Cell(1,4) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "79÷5=15, 4"

$cell = $t.Cell(1,5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "79÷6=13, 1") {
    Write-Host "WARNING: Cell(1,5) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "42÷4=10, 2"

$cell = $t.Cell(5,1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "89÷8=11, 1") {
    Write-Host "WARNING: Cell(5,1) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "81÷6=13, 3"

$cell = $t.Cell(5,2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "76÷8=9, 4") {
    Write-Host "WARNING: Cell(5,2) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "76÷9=8, 4"

$cell = $t.Cell(5,3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "21÷5=4, 1") {
    Write-Host "WARNING: Cell(5,3) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "98÷2=49, 0"

$cell = $t.Cell(5,4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "52÷8=6, 4") {
    Write-Host "WARNING: Cell(5,4) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "33÷8=4, 1"

$cell = $t.Cell(5,5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "60÷8=7, 4") {
    Write-Host "WARNING: Cell(5,5) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "77÷5=15, 2"

$cell = $t.Cell(9,1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "17÷9=1, 8") {
    Write-Host "WARNING: Cell(9,1) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "71÷8=8, 7"

$cell = $t.Cell(9,2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "66÷2=33, 0") {
    Write-Host "WARNING: Cell(9,2) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "36÷8=4, 4"

$cell = $t.Cell(9,3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "75÷2=37, 1") {
    Write-Host "WARNING: Cell(9,3) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "91÷7=13, 0"

$cell = $t.Cell(9,4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "88÷9=9, 7") {
    Write-Host "WARNING: Cell(9,4) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "55÷5=11, 0"

$cell = $t.Cell(9,5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "68÷6=11, 2") {
    Write-Host "WARNING: Cell(9,5) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "20÷9=2, 2"

$cell = $t.Cell(13,1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "70÷7=10, 0") {
    Write-Host "WARNING: Cell(13,1) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "91÷7=13, 0"

$cell = $t.Cell(13,2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "85÷9=9, 4") {
    Write-Host "WARNING: Cell(13,2) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "61÷5=12, 1"

$cell = $t.Cell(13,3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "46÷9=5, 1") {
    Write-Host "WARNING: Cell(13,3) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "59÷4=14, 3"

$cell = $t.Cell(13,4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "77÷5=15, 2") {
    Write-Host "WARNING: Cell(13,4) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "62÷5=12, 2"

$cell = $t.Cell(13,5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "52÷6=8, 4") {
    Write-Host "WARNING: Cell(13,5) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "50÷2=25, 0"

$cell = $t.Cell(17,1)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "28÷8=3, 4") {
    Write-Host "WARNING: Cell(17,1) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "27÷6=4, 3"

$cell = $t.Cell(17,2)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "66÷7=9, 3") {
    Write-Host "WARNING: Cell(17,2) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "99÷3=33, 0"

$cell = $t.Cell(17,3)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "55÷2=27, 1") {
    Write-Host "WARNING: Cell(17,3) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "14÷6=2, 2"

$cell = $t.Cell(17,4)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "94÷9=10, 4") {
    Write-Host "WARNING: Cell(17,4) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "71÷5=14, 1"

$cell = $t.Cell(17,5)
if ($cell.Range.Text.TrimEnd([char]7,[char]13) -ne "65÷3=21, 2") {
    Write-Host "WARNING: Cell(17,5) unexpected text: $($cell.Range.Text)"
}
$cell.Range.Text = "32÷5=6, 2"

